$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, pushing existing rows 100-157 down to 101-158
$ws.Rows.Item(100).Insert(4)   # xlShiftDown = -4121, but any shift-down works for a whole row

# Populate the newly inserted row 100 with the new record
$ws.Cells.Item(100, 1).Value = 5
$ws.Cells.Item(100, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(100, 3).Value = "Maule"

$ws.Cells.Item(100, 4).Value = 44488
$ws.Cells.Item(100, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat

$ws.Cells.Item(100, 5).Value = 7
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100108
$ws.Cells.Item(100, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(100, 9).Value = 100108005
$ws.Cells.Item(100, 10).Value = "Piña"
$ws.Cells.Item(100, 11).Value = "Caramelo"
$ws.Cells.Item(100, 12).Value = "Segunda"
$ws.Cells.Item(100, 13).Value = 540
$ws.Cells.Item(100, 14).Value = 18000
$ws.Cells.Item(100, 15).Value = 18000
$ws.Cells.Item(100, 16).Value = 18000
$ws.Cells.Item(100, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(100, 18).Value = "Ecuador"
$ws.Cells.Item(100, 19).Value = 1286
$ws.Cells.Item(100, 20).Value = 14
